# Insert a new daily-price row at row 439 (Vega Monumental Concepción - Plátano),
# which pushes the existing rows 439-527 down to 440-528.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(439).Insert()

# Populate the newly inserted row 439 with the new daily record.
$ws.Range("A439").Value = 11
$ws.Range("B439").Value = "Vega Monumental Concepción"
$ws.Range("C439").Value = "Bíobío"
$ws.Range("D439").Value = 44754
$ws.Range("E439").Value = 8
$ws.Range("F439").Value = "Fruta"
$ws.Range("G439").Value = 100108
$ws.Range("H439").Value = "Tropicales y subtropicales"
$ws.Range("I439").Value = 100108006
$ws.Range("J439").Value = "Plátano"
$ws.Range("K439").Value = "Sin especificar"
$ws.Range("L439").Value = "Pintón"
$ws.Range("M439").Value = 1000
$ws.Range("N439").Value = 22000
$ws.Range("O439").Value = 23000
$ws.Range("P439").Value = 22500
$ws.Range("Q439").Value = "$/caja 20 kilos"
$ws.Range("R439").Value = "Ecuador"
$ws.Range("S439").Value = 1125
$ws.Range("T439").Value = 20
